$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns I and J - copy H1's formatting (bold, centered, bordered header style)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2..26: I = 1 (with H25 -> I25 = 4 special case), J = H value (H25 -> J25 = 5 special case)
for ($r = 2; $r -le 26; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2

    if ($r -eq 25) {
        $ws.Cells.Item($r, 9).Value = 4
        $ws.Cells.Item($r, 10).Value = 5
    } else {
        $ws.Cells.Item($r, 9).Value = 1
        $ws.Cells.Item($r, 10).Value = $hVal
    }
}
